$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2" = 0.556007222541788
    "C2" = 0.612606004275779
    "K2" = 0.438282499020205
    "L2" = 0.595756191953926
    "N2" = 0.508231644030169
    "B3" = 0.50312291438834
    "K3" = 0.351050882580874
    "L3" = 0.604399737467109
    "N3" = 0.447300643788012
    "B4" = 0.609601586795904
    "K4" = 0.571096670838126
    "L4" = 0.778963825426238
    "N4" = 0.472900191628792
    "B5" = 0.347853243460036
    "N5" = 0.306445646731996
    "B6" = 0.410626908494325
    "K6" = 0.319383802321488
    "L6" = 0.389233362357354
    "N6" = 0.40055514051731
}

foreach ($cellRef in $values.Keys) {
    $ws.Range($cellRef).Value = $values[$cellRef]
}
